$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format column D as Text so numeric-looking price strings (e.g. "322.44")
# are stored as literal text, matching the source inlineStr cells, not coerced
# into floating point numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "47.875.96"
$ws.Range("E2").Value = "  +0.63%  "
$ws.Range("D3").Value = "2.503.44"
$ws.Range("E3").Value = "  +0.49%  "
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").Value = "322.44"
$ws.Range("E5").Value = "  +0.03%  "
$ws.Range("D6").Value = "109.49"
$ws.Range("E6").Value = "  +4.28%  "
$ws.Range("D7").Value = "0.524"
$ws.Range("E7").Value = "  -0.23%  "
$ws.Range("D8").Value = "0.999"
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("D9").Value = "0.544"
$ws.Range("E9").Value = "  +0.30%  "
$ws.Range("D10").Value = "39.59"
$ws.Range("E10").Value = "  +3.94%  "
$ws.Range("D11").Value = "0.0812"
$ws.Range("E11").Value = "  -0.19%  "
$ws.Range("E12").Value = "  +0.69%  "
$ws.Range("D13").Value = "18.65"
$ws.Range("E13").Value = "  +1.84%  "
$ws.Range("D14").Value = "7.23"
$ws.Range("E14").Value = "  +0.95%  "
$ws.Range("D15").Value = "2.887.04"
$ws.Range("E15").Value = "  +0.16%  "
$ws.Range("D16").Value = "2.504.69"
$ws.Range("E16").Value = "  +0.41%  "
$ws.Range("D17").Value = "0.849"
$ws.Range("E17").Value = "  +0.13%  "
$ws.Range("D18").Value = "47.767.32"
$ws.Range("E18").Value = "  +0.66%  "
$ws.Range("D19").Value = "13.34"
$ws.Range("E19").Value = "  +4.21%  "
$ws.Range("D20").Value = "6.66"
$ws.Range("E20").Value = "  +0.96%  "
$ws.Range("D21").Value = "0.0₃0944"
$ws.Range("E21").Value = "  +0.60%  "
$ws.Range("D22").Value = "2.76"
$ws.Range("E22").Value = "  +15.19%  "
$ws.Range("D23").Value = "70.77"
$ws.Range("E23").Value = "  +0.10%  "
$ws.Range("D24").Value = "247.67"
$ws.Range("E24").Value = "  -1.56%  "
$ws.Range("E25").Value = "  -0.49%  "
$ws.Range("E26").Value = "  +0.06%  "
$ws.Range("D27").Value = "25.86"
$ws.Range("E27").Value = "  -1.25%  "
$ws.Range("B28").Value = "Cosmos"
$ws.Range("C28").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D28").Value = "10.02"
$ws.Range("E28").Value = "  -0.14%  "
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").Value = "2.20"
$ws.Range("E29").Value = "  -3.05%  "
$ws.Range("D30").Value = "0.139"
$ws.Range("E30").Value = "  +2.83%  "
$ws.Range("D31").Value = "34.91"
$ws.Range("E31").Value = "  -0.77%  "
$ws.Range("D32").Value = "49.97"
$ws.Range("E32").Value = "  +1.14%  "
$ws.Range("D33").Value = "20.27"
$ws.Range("E33").Value = "  +2.54%  "
$ws.Range("D34").Value = "5.35"
$ws.Range("E34").Value = "  -0.25%  "
$ws.Range("D35").Value = "0.0790"
$ws.Range("E35").Value = "  +0.73%  "
$ws.Range("E36").Value = "  +0.00%  "
$ws.Range("D37").Value = "4.74"
$ws.Range("E37").Value = "  +2.20%  "
$ws.Range("D38").Value = "1.97"
$ws.Range("E38").Value = "  -0.15%  "
$ws.Range("D39").Value = "2.96"
$ws.Range("E39").Value = "  -1.22%  "
$ws.Range("E40").Value = "  +0.32%  "
$ws.Range("B41").Value = "EnergySwap"
$ws.Range("C41").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D41").Value = "22.37"
$ws.Range("E41").Value = "  +5.05%  "
$ws.Range("B42").Value = "WEMIXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D42").Value = "2.22"
$ws.Range("E42").Value = "  -1.41%  "
$ws.Range("D43").Value = "119.85"
$ws.Range("E43").Value = "  -1.77%  "
$ws.Range("E44").Value = "  +0.50%  "
$ws.Range("D45").Value = "1.996.75"
$ws.Range("E45").Value = "  +1.47%  "
$ws.Range("D46").Value = "3.05"
$ws.Range("E46").Value = "  +2.31%  "
$ws.Range("E47").Value = "  -2.99%  "
$ws.Range("D48").Value = "1.81"
$ws.Range("E48").Value = "  -0.33%  "
$ws.Range("E49").Value = "  -1.64%  "
$ws.Range("D50").Value = "5.24"
$ws.Range("E50").Value = "  -0.66%  "
$ws.Range("D51").Value = "56.77"
$ws.Range("E51").Value = "  +3.48%  "

# Column D cells were explicitly given a "@" (Text) number format above so that
# values such as "322.44" persist as text instead of being auto-coerced to
# numbers. Re-apply the original (unstyled) cell formatting on top so the saved
# file keeps the same default style as every other data cell in the sheet -
# only the text VALUE should differ from the source workbook, not the format.
$ws.Range("B2:B51").Copy() | Out-Null
$ws.Range("D2:D51").PasteSpecial(-4122)
$excel.CutCopyMode = 0
